$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 210 (shifts existing rows 210..229 down to 211..230)
$ws.Rows.Item(210).Insert()

# Populate the newly inserted row 210 with the new data record
$ws.Cells.Item(210, 1).Value = 7
$ws.Cells.Item(210, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(210, 3).Value = "Ñuble"
$ws.Cells.Item(210, 4).Value = 44769
$ws.Cells.Item(210, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(210, 5).Value = 16
$ws.Cells.Item(210, 6).Value = "Fruta"
$ws.Cells.Item(210, 7).Value = 100108
$ws.Cells.Item(210, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(210, 9).Value = 100108005
$ws.Cells.Item(210, 10).Value = "Piña"
$ws.Cells.Item(210, 11).Value = "Caramelo"
$ws.Cells.Item(210, 12).Value = "Segunda"
$ws.Cells.Item(210, 13).Value = 100
$ws.Cells.Item(210, 14).Value = 19000
$ws.Cells.Item(210, 15).Value = 20000
$ws.Cells.Item(210, 16).Value = 19500
$ws.Cells.Item(210, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(210, 18).Value = "Ecuador"
$ws.Cells.Item(210, 19).Value = 1393
$ws.Cells.Item(210, 20).Value = 14
